# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorganizes the "Estado de Cuenta" detail table (rows 16-33) so the rows
# are grouped by mora period (2111, then 2112, then 2201) instead of being
# grouped by worker. The same worker/period/value combinations are kept -
# only their row positions (and a couple of values that move along with
# them) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows     = @(16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33)

$docNums  = @(
    "33102996", "45460282", "87880053", "1047385464", "1002244348", "30764183",
    "33102996", "45460282", "87880053", "1047385464", "1002244348", "30764183",
    "33102996", "45460282", "87880053", "1047385464", "1002244348", "30764183"
)

$names    = @(
    "SHIRLIS ALVAREZ ESCALANTE", "ALIEIDA DEL ROSARIO JULIO PINILLA", "JUSTO TELLO",
    "YOCELIN GOMEZ BARRIOS", "FEDERICO DE JESUS DIAZ CASTRO", "GLORIA PATRICIA CASTRO CARRILLO",
    "SHIRLIS ALVAREZ ESCALANTE", "ALIEIDA DEL ROSARIO JULIO PINILLA", "JUSTO TELLO",
    "YOCELIN GOMEZ BARRIOS", "FEDERICO DE JESUS DIAZ CASTRO", "GLORIA PATRICIA CASTRO CARRILLO",
    "SHIRLIS ALVAREZ ESCALANTE", "ALIEIDA DEL ROSARIO JULIO PINILLA", "JUSTO TELLO",
    "YOCELIN GOMEZ BARRIOS", "FEDERICO DE JESUS DIAZ CASTRO", "GLORIA PATRICIA CASTRO CARRILLO"
)

$periods  = @(
    "2111", "2111", "2111", "2111", "2111", "2111",
    "2112", "2112", "2112", "2112", "2112", "2112",
    "2201", "2201", "2201", "2201", "2201", "2201"
)

$moras    = @(
    36341, 36341, 120000, 52000, 52000, 52000,
    36341, 36341, 120000, 52000, 52000, 52000,
    31495, 24227, 104000, 34666, 34666, 34666
)

$salarios = @(
    908526, 908526, 3000000, 1300000, 1300000, 1300000,
    908526, 908526, 3000000, 1300000, 1300000, 1300000,
    908526, 908526, 3000000, 1300000, 1300000, 1300000
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 3).Value = $docNums[$i]
    $ws.Cells.Item($r, 4).Value = $names[$i]
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = $moras[$i]
    $ws.Cells.Item($r, 7).Value = $salarios[$i]
}
